# Adds a new "TC005__ValidateUpdateAPILatency" data block to the sheet,
# mirroring the existing TC blocks (header row + name/job row + values row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the first test-case block (rows 1-3) down to
# rows 13-15 (row 12 stays blank, matching the existing blank-row separator
# pattern used between the other blocks).
$ws.Range("A1").Copy($ws.Range("A13"))
$ws.Range("A2:B3").Copy($ws.Range("A14"))

# Replace the header text on the new block with the new test case name.
$ws.Range("A13").Value2 = "TC005__ValidateUpdateAPILatency"

# Match the saved selection/active cell.
$ws.Range("A13").Select()
